# Convertation to the sheet
# Convert the numeric header/time-slot placeholders to actual weekday names
# and class time ranges, matching the real MEC-3A timetable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): weekday names
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# Column A: time slots for each class period
$ws.Range("A2").Value = "7:00"
$ws.Range("A3").Value = "7:50"
$ws.Range("A4").Value = "8:40"
$ws.Range("A5").Value = "9:30"
$ws.Range("A6").Value = "10:40"
$ws.Range("A7").Value = "11:30"
$ws.Range("A8").Value = "13:00"
$ws.Range("A9").Value = "13:50"
$ws.Range("A10").Value = "14:40"
$ws.Range("A11").Value = "15:30"
$ws.Range("A12").Value = "16:40"
$ws.Range("A13").Value = "17:30"

# Move class entries to the correct (newly re-labeled) slots
$ws.Range("D3").Value = "EAP"

$ws.Range("E4").Value = "Circuitos Elétricos 2"

$ws.Range("C5").Value = "-"

$ws.Range("B9").Value = "-"
$ws.Range("C9").Value = "Desenho Técnico"

$ws.Range("E10").Value = "EAP"

$ws.Range("C13").Value = "-"
$ws.Range("E13").Value = "-"
